$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 01:05"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 2).Value = 1618820
$ws.Cells.Item(4, 3).Value = 26097
$ws.Cells.Item(4, 4).Value = 381782
$ws.Cells.Item(4, 5).Value = 1140812
$ws.Cells.Item(4, 7).Value = 1290
$ws.Cells.Item(4, 8).Value = 96226

# Row 6: Brasil -> Brasil
$ws.Cells.Item(6, 2).Value = 310087
$ws.Cells.Item(6, 3).Value = 16730
$ws.Cells.Item(6, 4).Value = 125960
$ws.Cells.Item(6, 5).Value = 164080
$ws.Cells.Item(6, 7).Value = 1153
$ws.Cells.Item(6, 8).Value = 20047

# Row 51: Argentina -> Argentina
$ws.Cells.Item(51, 2).Value = 9931
$ws.Cells.Item(51, 3).Value = 648
$ws.Cells.Item(51, 5).Value = 6483
$ws.Cells.Item(51, 7).Value = 13
$ws.Cells.Item(51, 8).Value = 416

# Row 52: Chequia -> Chequia
$ws.Cells.Item(52, 2).Value = 8754
$ws.Cells.Item(52, 3).Value = 33
$ws.Cells.Item(52, 4).Value = 5926
$ws.Cells.Item(52, 5).Value = 2522

# Row 54: Noruega -> Noruega
$ws.Cells.Item(54, 2).Value = 8309
$ws.Cells.Item(54, 3).Value = 28
$ws.Cells.Item(54, 5).Value = 8042

# Row 79: Senegal -> Senegal
$ws.Cells.Item(79, 5).Value = 1528
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 33

# Row 110: Mali -> Mali
$ws.Cells.Item(110, 2).Value = 947
$ws.Cells.Item(110, 3).Value = 16
$ws.Cells.Item(110, 4).Value = 558
$ws.Cells.Item(110, 5).Value = 329
$ws.Cells.Item(110, 7).Value = 5
$ws.Cells.Item(110, 8).Value = 60

# Row 115: Zambia -> Venezuela
$ws.Cells.Item(115, 1).Value = "Venezuela"
$ws.Cells.Item(115, 2).Value = 882
$ws.Cells.Item(115, 3).Value = 58
$ws.Cells.Item(115, 4).Value = 262
$ws.Cells.Item(115, 5).Value = 610
$ws.Cells.Item(115, 8).Value = 10

# Row 116: Paraguay -> Zambia
$ws.Cells.Item(116, 1).Value = "Zambia"
$ws.Cells.Item(116, 2).Value = 866
$ws.Cells.Item(116, 3).Value = 34
$ws.Cells.Item(116, 4).Value = 302
$ws.Cells.Item(116, 5).Value = 557
$ws.Cells.Item(116, 8).Value = 7

# Row 117: Venezuela -> Paraguay
$ws.Cells.Item(117, 1).Value = "Paraguay"
$ws.Cells.Item(117, 2).Value = 836
$ws.Cells.Item(117, 3).Value = 3
$ws.Cells.Item(117, 4).Value = 256
$ws.Cells.Item(117, 5).Value = 569
$ws.Cells.Item(117, 8).Value = 11

# Row 149: Liberia -> Guayana Francesa
$ws.Cells.Item(149, 1).Value = "Guayana Francesa"
$ws.Cells.Item(149, 2).Value = 249
$ws.Cells.Item(149, 3).Value = 12
$ws.Cells.Item(149, 4).Value = 137
$ws.Cells.Item(149, 5).Value = 111
$ws.Cells.Item(149, 8).Value = 1

# Row 150: Guayana Francesa -> Liberia
$ws.Cells.Item(150, 1).Value = "Liberia"
$ws.Cells.Item(150, 2).Value = 240
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(150, 4).Value = 131
$ws.Cells.Item(150, 5).Value = 86
$ws.Cells.Item(150, 8).Value = 23

# Row 167: Trinidad yTobago -> Islas Caimanes
$ws.Cells.Item(167, 1).Value = "Islas Caimanes"
$ws.Cells.Item(167, 2).Value = 121
$ws.Cells.Item(167, 3).Value = 10
$ws.Cells.Item(167, 4).Value = 55
$ws.Cells.Item(167, 5).Value = 65
$ws.Cells.Item(167, 8).Value = 1

# Row 168: Islas Caimanes -> Trinidad yTobago
$ws.Cells.Item(168, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(168, 2).Value = 116
$ws.Cells.Item(168, 4).Value = 107
$ws.Cells.Item(168, 5).Value = 1
$ws.Cells.Item(168, 8).Value = 8

# Row 209: Groenlandia -> Montserrat
$ws.Cells.Item(209, 1).Value = "Montserrat"
$ws.Cells.Item(209, 4).Value = 10
$ws.Cells.Item(209, 8).Value = 1

# Row 210: Seychelles -> Groenlandia
$ws.Cells.Item(210, 1).Value = "Groenlandia"

# Row 211: Montserrat -> Seychelles
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Row 214: Bonaire, San Eustaquio y Saba -> Sahara Occidental
$ws.Cells.Item(214, 1).Value = "Sahara Occidental"

# Row 215: Sahara Occidental -> Bonaire, San Eustaquio y Saba
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
